$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-order worksheets.
#    Before: AddOpportunity, AppName, ModuleName, Users, AddContact,
#            OppDealTeamMembers, EngDealTeamMembers, OverlimitMessage
#    After:  Users, AppName, ModuleName, AddOpportunity, AddContact,
#            OppDealTeamMembers, EngDealTeamMembers, OverlimitMessage
# ---------------------------------------------------------------------------

# Move "AddOpportunity" so it sits immediately before "AddContact"
# (i.e. right after "ModuleName", where "Users" used to be).
$addOpportunity = $wb.Worksheets.Item("AddOpportunity")
$addContact = $wb.Worksheets.Item("AddContact")
[void]$addOpportunity.Move($addContact)

# Move "Users" to be the very first sheet in the workbook.
$users = $wb.Worksheets.Item("Users")
$firstSheet = $wb.Worksheets.Item(1)
[void]$users.Move($firstSheet)

# ---------------------------------------------------------------------------
# 2. Data edits.
#    New shared strings are appended in first-use order, so the cell writes
#    below are ordered to reproduce the target sharedStrings.xml layout
#    (Ryan Mahlan, Rebecca Hu, Timothy Kang, Alan Test).
# ---------------------------------------------------------------------------

# EngDealTeamMembers: swap the team member name.
$engDealTeamWs = $wb.Worksheets.Item("EngDealTeamMembers")
$engDealTeamWs.Range("A2").Value = "Ryan Mahlan"

# OppDealTeamMembers: append two new deal team members.
$oppDealTeamWs = $wb.Worksheets.Item("OppDealTeamMembers")
$oppDealTeamWs.Range("A29").Value = "Rebecca Hu"
$oppDealTeamWs.Range("A30").Value = "Timothy Kang"
$oppDealTeamWs.Range("A30").WrapText = $true

# AddContact: swap the contact name.
$addContactWs = $wb.Worksheets.Item("AddContact")
$addContactWs.Range("A2").Value = "Alan Test"

# ---------------------------------------------------------------------------
# 3. Restore per-sheet cursor/selection state.
#    Selecting a range activates that sheet, so sheets are touched in an
#    order that leaves "AddContact" as the final active tab (matching the
#    target workbook's activeTab / tabSelected state).
# ---------------------------------------------------------------------------

$usersWs = $wb.Worksheets.Item("Users")
[void]$usersWs.Range("F16").Select()

$oppDealTeamWs = $wb.Worksheets.Item("OppDealTeamMembers")
[void]$oppDealTeamWs.Range("B29").Select()

$engDealTeamWs = $wb.Worksheets.Item("EngDealTeamMembers")
[void]$engDealTeamWs.Range("H23").Select()

$addContactWs = $wb.Worksheets.Item("AddContact")
[void]$addContactWs.Activate()
[void]$addContactWs.Range("C6").Select()
